# Insert a new data row for "San Juan, Puerto Rico" (SJU) right before the
# existing "Amman, Jordan" row (row 225), shifting every subsequent row down
# by one. This also grows the used range from A1:G319 to A1:G320.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 225; everything from old row 225 onward
# (Amman ... Suva) slides down to 226 ... 320.
$ws.Rows.Item(225).Insert()

# Copy the (now shifted-down) formatting of what used to be row 225 so the
# new row matches the rest of the table's look (bordered / bold colo code
# in column A, etc.).
$ws.Range("A226").Copy()
$ws.Range("A225").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new row's values.
$ws.Range("A225").Value2 = "SJU"
$ws.Range("B225").Value2 = "San Juan, Puerto Rico"
$ws.Range("C225").Value2 = 18.411391
$ws.Range("D225").Value2 = -66.10279300000001
$ws.Range("E225").Value2 = "PR"
$ws.Range("F225").Value2 = "North America"
$ws.Range("G225").Value2 = "San Juan"
